# Insert a new weekly price record as row 305 on the single data sheet,
# pushing all subsequent rows (old 305-366) down by one (new 306-367).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 305; Excel shifts rows 305:366 down to 306:367
# and the used-range dimension grows from R366 to R367 automatically.
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A305").Value = 9
$ws.Range("B305").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C305").Value = "Metropolitana"
$ws.Range("D305").Value = 44932
$ws.Range("E305").Value = 13
$ws.Range("F305").Value = 100112021
$ws.Range("G305").Value = "Ají"
$ws.Range("H305").Value = "Americana (o)"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 70
$ws.Range("K305").Value = 21000
$ws.Range("L305").Value = 23000
$ws.Range("M305").Value = 22000
$ws.Range("N305").Value = "$/caja 15 kilos"
$ws.Range("O305").Value = "Región Metropolitana"
$ws.Range("P305").Value = 1467
$ws.Range("Q305").Value = 15
$ws.Range("R305").Value = "Hortaliza"
